# Append two new daily rows (2025-11-24 and 2025-11-25, Excel serials
# 45985 / 45986) to the bottom of each of the six worksheets' A:B data
# tables (date / remn_amt), matching the existing row 115 -> row 116/117
# layout and formatting (date column uses the existing "YYYY-MM-DD HH:MM:SS"
# number format already applied to the rest of column A).

$wb = $excel.ActiveWorkbook

# sheet index => [row116_B, row117_B]; row116/117 A values are the same
# two dates (45985, 45986) on every sheet.
$newData = @{
    1 = @(1287651, 1280640)   # 현대차
    2 = @(915695,  915802)    # 기아
    3 = @(424331,  424366)    # 현대모비스
    4 = @(900126,  897401)    # 삼성중공업
    5 = @(3151079, 3088022)   # HD현대중공업
    6 = @(865704,  839178)    # 한화오션
}

foreach ($idx in 1..6) {
    $ws = $wb.Worksheets.Item($idx)
    $vals = $newData[$idx]

    $ws.Range("A116").Value() = 45985
    $ws.Range("B116").Value() = $vals[0]

    $ws.Range("A117").Value() = 45986
    $ws.Range("B117").Value() = $vals[1]

    # Match the date formatting already used for the rest of column A.
    $ws.Range("A116:A117").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
